# Small analysis tools - sector coupled network
# Insert 3 new carrier rows (H2, H2 electrolysis, H2 fuel cell) into the
# "Carrier" sheet just above the existing "diesel" row, and make the
# "Carrier" sheet the active/selected sheet (instead of "Store").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Carrier")

# Insert three blank rows at 14, pushing the existing "diesel" row (and
# anything below it) down to row 17.
$ws.Rows("14:16").Insert()
$ws.Rows("14:16").RowHeight = 18.75

$ws.Range("A14").Value = "H2"
$ws.Range("B14").Value = "#16942d"

$ws.Range("A15").Value = "H2 electrolysis"
$ws.Range("B15").Value = "#16942d"

$ws.Range("A16").Value = "H2 fuel cell"
$ws.Range("B16").Value = "#16942d"

# Make the "Carrier" sheet the active tab (was "Store").
$ws.Activate()
